# Generate Report for Handback
#
# Updates the handoff/handback timestamp columns for the
# "509d6e4f-e05b-401f-b527-53db60051336.md" row (row 3 on every sheet)
# across the Overview / zh-cn / de-de worksheets, reflecting a freshly
# generated handback report.

$wb = $excel.ActiveWorkbook

# --- Overview sheet --------------------------------------------------
# Column G = "Latest HO Xliff Generate Date"
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G3").Value = "2016-11-15 16:55:01"

# --- zh-cn sheet -------------------------------------------------------
# Column H = "Correspond Handoff Datetime", Column K = "Correspond Handback DateTime"
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H3").Value = "2016-11-15 16:54:45"
$wsZhCn.Range("K3").Value = "2016-11-15 16:55:44"

# --- de-de sheet -------------------------------------------------------
# Column H = "Correspond Handoff Datetime", Column K = "Correspond Handback DateTime"
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H3").Value = "2016-11-15 16:55:01"
$wsDeDe.Range("K3").Value = "2016-11-15 16:56:04"
